$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Prato do chefe - teste número 05", "Arroz Parboilizado", 150, "Un"),
    @("Prato do chefe - teste número 05", "Alho Poro", 50, "Un"),
    @("Prato do chefe - teste número 05", "Bife do Vazio", 300, "g"),
    @("Prato do chefe - teste número 05", "Batata frita corte fino", 148, "Kg"),
    @("Prato do chefe - teste número 05", "Cebola Roxa", 101, "g")
)

$startRow = 10
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}
